# Adds non-road hydrogen vehicles, adjusts biofuels share, H2 shares to match

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsRHPF  = $wb.Worksheets.Item("RHPF")

# Update formulas on the RHPF sheet:
#   B2 (electrolysis column, electrolysis row) and B8 (electrolysis column,
#   "natural gas reforming with CCS" row) now pull a computed share, and
#   B7 (electrolysis column, "electrolysis with guaranteed clean electricity"
#   row) is adjusted so the column still sums to 1.
$wsRHPF.Range("B2").Formula = "=0.051/0.2092/0.795"
$wsRHPF.Range("B8").Formula = "=0.051/0.2092/0.795"
$wsRHPF.Range("B7").Formula = "=1-B2-B8"

# Make the RHPF sheet the active tab/sheet (previously "About" was active),
# with the selection left on B8.
$wsRHPF.Activate()
$wsRHPF.Range("B8").Select()
